$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update F column (time_taken) timestamps in "data" sheet
$ws.Range("F2").Value = "2021-10-05 14:33:50.382611"
$ws.Range("F3").Value = "2021-10-05 14:33:50.382621"
$ws.Range("F4").Value = "2021-10-05 14:33:50.382624"
$ws.Range("F5").Value = "2021-10-05 14:33:50.382627"
$ws.Range("F6").Value = "2021-10-05 14:33:50.382630"
$ws.Range("F7").Value = "2021-10-05 14:33:50.382633"
$ws.Range("F8").Value = "2021-10-05 14:33:50.382636"
$ws.Range("F9").Value = "2021-10-05 14:33:50.382638"
$ws.Range("F10").Value = "2021-10-05 14:33:50.382641"
$ws.Range("F11").Value = "2021-10-05 14:33:50.382644"
$ws.Range("F12").Value = "2021-10-05 14:33:50.382647"
$ws.Range("F13").Value = "2021-10-05 14:33:50.382650"
$ws.Range("F14").Value = "2021-10-05 14:33:50.382653"
$ws.Range("F15").Value = "2021-10-05 14:33:50.382655"
$ws.Range("F16").Value = "2021-10-05 14:33:50.382658"
$ws.Range("F17").Value = "2021-10-05 14:33:50.382661"
$ws.Range("F18").Value = "2021-10-05 14:33:50.382664"
$ws.Range("F19").Value = "2021-10-05 14:33:50.382666"
$ws.Range("F20").Value = "2021-10-05 14:33:50.382669"

# Add a new "metadata" worksheet, positioned after "data"
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1), bold/bordered like the "data" sheet header
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row 2
$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Eye Anterior Segment Abnormalities"
$meta.Range("C2").Value = 43
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"
$meta.Range("D2").ClearFormats()
$meta.Range("E2").Value = "2020-12-29T07:36:00.040753Z"
$meta.Range("F2").Value = "2021-10-05 14:33:50.379135"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/43/?format=json"

# Restore "data" as the active sheet (matches original workbook view state)
$ws.Activate()
